# Auto-generated edit script
# Applies the "one day passed" update to the 剩余(remaining) and 开始时间(start date)
# columns on the active worksheet, per commit message:
# 自动更新Excel文件 - 2025-12-29 23:13:19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("E2").Value = 13
    $ws.Range("E3").Value = 13
    $ws.Range("E4").Value = 13
    $ws.Range("E5").Value = 5
    $ws.Range("E6").Value = 13
    $ws.Range("E7").Value = 5
    $ws.Range("E8").Value = 13
    $ws.Range("E9").Value = 5
    $ws.Range("E10").Value = 6
    $ws.Range("E11").Value = 13
    $ws.Range("E12").Value = 5
    $ws.Range("E13").Value = 13
    $ws.Range("E14").Value = 13
    $ws.Range("E15").Value = 13
    $ws.Range("E16").Value = 9
    $ws.Range("E17").Value = 5
    $ws.Range("E18").Value = 8
    $ws.Range("E19").Value = 8
    $ws.Range("E20").Value = 8
    $ws.Range("E21").Value = 8
    $ws.Range("E22").Value = 5
    $ws.Range("E23").Value = 5
    $ws.Range("E24").Value = 5
    $ws.Range("E25").Value = 5
    $ws.Range("E26").Value = 5
    $ws.Range("E27").Value = 7
    $ws.Range("F27").Value = 20251230
    $ws.Range("E28").Value = 8
    $ws.Range("E29").Value = 8
    $ws.Range("E30").Value = 8
    $ws.Range("E31").Value = 8
    $ws.Range("E32").Value = 8
    $ws.Range("E33").Value = 8
    $ws.Range("E34").Value = 8
    $ws.Range("E35").Value = 8
    $ws.Range("E37").Value = 8
    $ws.Range("E38").Value = 8
    $ws.Range("E39").Value = 8
    $ws.Range("E40").Value = 6
    $ws.Range("E41").Value = 6
    $ws.Range("E42").Value = 8
    $ws.Range("E43").Value = 5
    $ws.Range("E44").Value = 6
    $ws.Range("E45").Value = 5
    $ws.Range("E46").Value = 6
    $ws.Range("E47").Value = 8
    $ws.Range("E48").Value = 6
    $ws.Range("E49").Value = 7
    $ws.Range("F49").Value = 20251230
    $ws.Range("E50").Value = 3
    $ws.Range("E51").Value = 3
    $ws.Range("E52").Value = 3
    $ws.Range("E53").Value = 3
    $ws.Range("E54").Value = 3
    $ws.Range("E55").Value = 3
    $ws.Range("E56").Value = 3
    $ws.Range("E57").Value = 3
    $ws.Range("E58").Value = 7
    $ws.Range("E59").Value = 7
    $ws.Range("E60").Value = 7
    $ws.Range("E61").Value = 7
    $ws.Range("F61").Value = 20251230
    $ws.Range("E62").Value = 7
    $ws.Range("E63").Value = 7
    $ws.Range("E64").Value = 7
    $ws.Range("E65").Value = 8
    $ws.Range("E66").Value = 8
    $ws.Range("E67").Value = 8
    $ws.Range("E68").Value = 8
    $ws.Range("E69").Value = 8
    $ws.Range("E70").Value = 9
    $ws.Range("E71").Value = 9
    $ws.Range("E72").Value = 9
    $ws.Range("E73").Value = 9
    $ws.Range("E74").Value = 9
    $ws.Range("E75").Value = 9
    $ws.Range("E76").Value = 9
    $ws.Range("E77").Value = 2
    $ws.Range("E78").Value = 2
    $ws.Range("E79").Value = 2
    $ws.Range("E80").Value = 2
    $ws.Range("E81").Value = 2
    $ws.Range("E82").Value = 2
    $ws.Range("E83").Value = 2
    $ws.Range("E84").Value = 2
    $ws.Range("E85").Value = 2
    $ws.Range("E86").Value = 2
    $ws.Range("E87").Value = 6
    $ws.Range("E88").Value = 6
    $ws.Range("E89").Value = 6
    $ws.Range("E90").Value = 6
    $ws.Range("E91").Value = 5
    $ws.Range("E92").Value = 6
    $ws.Range("E93").Value = 2
    $ws.Range("E94").Value = 2
    $ws.Range("E95").Value = 1
    $ws.Range("E96").Value = 9
    $ws.Range("E97").Value = 9
    $ws.Range("E98").Value = 9
    $ws.Range("E99").Value = 9
